$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")
$ws.Range("A1").Value = "test"
